# fix(gui) step 1 and 2
#
# Step 1: bump the date shown in A1 by one day.
# Step 2: update the unit prices in the two price tables
#         ("REFORZADA / ZINCADA" and "PINTADA BLANCA").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1 - date in A1 (was 45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2 - PRECIO C/U values
$ws.Range("D14").Value = 248.739
$ws.Range("D15").Value = 381.532
$ws.Range("D38").Value = 457.837
$ws.Range("D39").Value = 491.531
